$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 17928.947
$ws.Range("I74").Value = 18526.8
$ws.Range("J74").Value = 15687
$ws.Range("K74").Value = 18526.8
$ws.Range("L74").Value = 15687
$ws.Range("M74").Value = -17590.8
$ws.Range("N74").Value = -17559

$ws.Range("H77").Value = 17928.947
$ws.Range("I77").Value = 18526.8
$ws.Range("J77").Value = 15687
$ws.Range("K77").Value = 92634
$ws.Range("L77").Value = 78435
$ws.Range("M77").Value = -87954
$ws.Range("N77").Value = -87795

$ws.Range("H100").Value = 2502
$ws.Range("I100").Value = 1005
$ws.Range("J100").Value = 3999
$ws.Range("K100").Value = 1005
$ws.Range("L100").Value = 3999
$ws.Range("M100").Value = -464
$ws.Range("N100").Value = -5081

$ws.Range("H107").Value = 421.44446
$ws.Range("I107").Value = 528.3333
$ws.Range("K107").Value = 528.3333
$ws.Range("M107").Value = 1391.6667

$ws.Range("H135").Value = 881.93335
$ws.Range("I135").Value = 322.6154
$ws.Range("K135").Value = 2903.5386
$ws.Range("M135").Value = -368.5386000000003

$ws.Range("H137").Value = 5615.0835
$ws.Range("I137").Value = 1812.6666
$ws.Range("J137").Value = 9417.5
$ws.Range("K137").Value = 5437.9998
$ws.Range("L137").Value = 28252.5
$ws.Range("M137").Value = -2887.9998
$ws.Range("N137").Value = -33352.5

$ws.Range("H138").Value = 360758
$ws.Range("I138").Value = 4028.125
$ws.Range("J138").Value = 546877.9399999999
$ws.Range("K138").Value = 12084.375
$ws.Range("L138").Value = 1640633.82
$ws.Range("M138").Value = -6944.375
$ws.Range("N138").Value = -1650913.82

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3562.1233
$ws.Range("I32").Value = 3141.352
$ws.Range("K32").Value = 3141.352
$ws.Range("M32").Value = -2854.352

$ws.Range("H74").Value = 258845.9
$ws.Range("I74").Value = 371733.12
$ws.Range("K74").Value = 371733.12
$ws.Range("M74").Value = -370859.12

$ws.Range("H77").Value = 258845.9
$ws.Range("I77").Value = 371733.12
$ws.Range("K77").Value = 1858665.6
$ws.Range("M77").Value = -1854297.6

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").ClearContents()
$ws.Range("N109").Value = 0

$ws.Range("H110").Value = 4117.3438
$ws.Range("I110").Value = 3509.2415
$ws.Range("K110").Value = 3509.2415
$ws.Range("M110").Value = -1464.2415

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 10402676
$ws.Range("I105").Value = 558052.9399999999
$ws.Range("K105").Value = 558052.9399999999
$ws.Range("M105").Value = -556305.9399999999

$ws.Range("H108").Value = 67999
$ws.Range("I108").Value = 67999
$ws.Range("K108").Value = 67999
$ws.Range("M108").Value = -64159

$ws.Range("H134").Value = 2423
$ws.Range("I134").Value = 1674.9
$ws.Range("K134").Value = 5024.700000000001
$ws.Range("M134").Value = -2489.700000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3998.3142
$ws.Range("I31").Value = 3507.64
$ws.Range("K31").Value = 3507.64
$ws.Range("M31").Value = -3212.64

$ws.Range("H34").Value = 3998.3142
$ws.Range("I34").Value = 3507.64
$ws.Range("K34").Value = 3507.64
$ws.Range("M34").Value = -3305.64

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1869.6364
$ws.Range("I5").Value = 527
$ws.Range("K5").Value = 1581
$ws.Range("M5").Value = -1469

$ws.Range("H50").Value = 1098.5625
$ws.Range("I50").Value = 699.5
$ws.Range("K50").Value = 2098.5
$ws.Range("M50").Value = -1617.5

$ws.Range("H53").Value = 1098.5625
$ws.Range("I53").Value = 699.5
$ws.Range("K53").Value = 2098.5
$ws.Range("M53").Value = -1617.5

$ws.Range("H69").Value = 3500
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()

$ws.Range("H70").Value = 3377.25
$ws.Range("I70").Value = 498
$ws.Range("J70").Value = 4337
$ws.Range("K70").Value = 1494
$ws.Range("L70").Value = 13011
$ws.Range("M70").Value = -1179
$ws.Range("N70").Value = -13641

$ws.Range("H72").Value = 3500
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()

$ws.Range("H73").Value = 3377.25
$ws.Range("I73").Value = 498
$ws.Range("J73").Value = 4337
$ws.Range("K73").Value = 1494
$ws.Range("L73").Value = 13011
$ws.Range("M73").Value = -402
$ws.Range("N73").Value = -15195

$ws.Range("H81").Value = 3674.625
$ws.Range("J81").Value = 6499.25
$ws.Range("L81").Value = 19497.75
$ws.Range("N81").Value = -21743.75

$ws.Range("H84").Value = 3674.625
$ws.Range("J84").Value = 6499.25
$ws.Range("L84").Value = 58493.25
$ws.Range("N84").Value = -69725.25

$ws.Range("H135").Value = 1869.6364
$ws.Range("I135").Value = 527
$ws.Range("K135").Value = 4743
$ws.Range("M135").Value = -2208

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 23130.8
$ws.Range("I46").Value = 999.5
$ws.Range("J46").Value = 37885
$ws.Range("K46").Value = 999.5
$ws.Range("L46").Value = 37885
$ws.Range("M46").Value = -843.5
$ws.Range("N46").Value = -38197

$ws.Range("H57").Value = 53787.145
$ws.Range("J57").Value = 53787.145
$ws.Range("L57").Value = 53787.145
$ws.Range("N57").Value = -55427.145

$ws.Range("I80").Value = 166668320
$ws.Range("J80").Value = 6500
$ws.Range("K80").Value = 166668320
$ws.Range("L80").Value = 6500
$ws.Range("M80").Value = -166667322
$ws.Range("N80").Value = -8496

$ws.Range("I83").Value = 166668320
$ws.Range("J83").Value = 6500
$ws.Range("K83").Value = 833341600
$ws.Range("L83").Value = 32500
$ws.Range("M83").Value = -833336608
$ws.Range("N83").Value = -42484

$ws.Range("H107").Value = 1546.9166
$ws.Range("I107").Value = 1759.625
$ws.Range("J107").Value = 1121.5
$ws.Range("K107").Value = 1759.625
$ws.Range("L107").Value = 1121.5
$ws.Range("M107").Value = 160.375
$ws.Range("N107").Value = -4961.5

$ws.Range("H132").Value = 1978.9445
$ws.Range("I132").Value = 1440
$ws.Range("J132").Value = 3380.2
$ws.Range("K132").Value = 4320
$ws.Range("L132").Value = 10140.6
$ws.Range("M132").Value = -1790
$ws.Range("N132").Value = -15200.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1161.5807
$ws.Range("I61").Value = 1033.6333
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 1033.6333
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -831.6333
$ws.Range("N61").Value = -5404

$ws.Range("H93").Value = 796.7
$ws.Range("I93").Value = 733.625
$ws.Range("K93").Value = 733.625
$ws.Range("M93").Value = 514.375

$ws.Range("H100").Value = 3991.9
$ws.Range("I100").Value = 4002.7334
$ws.Range("K100").Value = 4002.7334
$ws.Range("M100").Value = -3461.7334

$ws.Range("H113").Value = 1161.5807
$ws.Range("I113").Value = 1033.6333
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 1033.6333
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = 1136.3667
$ws.Range("N113").Value = -9340

$ws.Range("H132").Value = 4107.9297
$ws.Range("I132").Value = 3230.8125
$ws.Range("K132").Value = 9692.4375
$ws.Range("M132").Value = -7162.4375

$ws.Range("H136").Value = 3948.162
$ws.Range("I136").Value = 3633
$ws.Range("K136").Value = 10899
$ws.Range("M136").Value = -8349

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 47620004
$ws.Range("I100").Value = 924.1818
$ws.Range("J100").Value = 100000990
$ws.Range("K100").Value = 1848.3636
$ws.Range("L100").Value = 200001980
$ws.Range("M100").Value = -1307.3636
$ws.Range("N100").Value = -200003062

$ws.Range("H126").Value = 1976.8636
$ws.Range("I126").Value = 1727.0555
$ws.Range("J126").Value = 3101
$ws.Range("K126").Value = 5181.166499999999
$ws.Range("L126").Value = 9303
$ws.Range("M126").Value = -2711.166499999999
$ws.Range("N126").Value = -14243

$ws.Range("H132").Value = 7248691
$ws.Range("I132").Value = 8774015
$ws.Range("K132").Value = 26322045
$ws.Range("M132").Value = -26319515

$ws.Range("H136").Value = 26318190
$ws.Range("I136").Value = 29412808
$ws.Range("K136").Value = 88238424
$ws.Range("M136").Value = -88235874

